$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.650.27"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "1.593.72"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "0.246"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").Value = "19.62"
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "1.816.69"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "1.598.86"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("E14").Value = "  -2.41%  "
$ws.Range("D15").Value = "0.523"
$ws.Range("E15").Value = "  -2.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "26.620.39"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "209.28"
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "6.68"
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("D22").Value = "4.24"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("D24").Value = "8.89"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("D25").Value = "146.58"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "7.14"
$ws.Range("E27").Value = "  -3.92%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "15.29"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("D33").Value = "0.693"
$ws.Range("E33").Value = "  -4.43%  "
$ws.Range("D35").Value = "1.291.53"
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  -5.50%  "
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("D42").Value = "5.34"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("D44").Value = "63.48"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").Value = "1.729.50"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.900"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.16%  "
$ws.Range("D47").Value = "89.73"
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("D49").Value = "0.0983"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.03%  "
